$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

# A new QR-scan record was prepended to the table: insert a fresh row
# above the existing row 2 (pushing the old rows 2-4 down to 3-5) and
# populate it with the scanned record's fields.
$ws.Rows.Item(2).Insert()

$newRow = 2
$ws.Cells.Item($newRow, 1).Value = "iowi2cez33m"
$ws.Cells.Item($newRow, 2).Value = "8n59cw7k"
$ws.Cells.Item($newRow, 3).Value = "B"
$ws.Cells.Item($newRow, 4).Value = "Phường Tăng Nhơn Phú, Thành phố Hồ Chí Minh, 71300, Việt Nam"
$ws.Cells.Item($newRow, 5).Value = "https://www.google.com/maps/search/?api=1&query=10.839061,106.792777"
$ws.Cells.Item($newRow, 6).Value = "2025-08-22T09:39:08.817Z"

# note / phone / branch / cccd / customerCode are empty strings (not
# blank cells) in the source data. A plain "" assignment clears the
# cell entirely in Excel, so force literal empty text via the
# leading-apostrophe (quote-prefix) trick, then drop the quote-prefix
# style it applies so the cell's formatting stays plain.
$ws.Cells.Item($newRow, 7).Value = "'"
$ws.Cells.Item($newRow, 7).Style = "Normal"
$ws.Cells.Item($newRow, 8).Value = "'"
$ws.Cells.Item($newRow, 8).Style = "Normal"
$ws.Cells.Item($newRow, 9).Value = "'"
$ws.Cells.Item($newRow, 9).Style = "Normal"
$ws.Cells.Item($newRow, 10).Value = "'"
$ws.Cells.Item($newRow, 10).Style = "Normal"
$ws.Cells.Item($newRow, 11).Value = "'"
$ws.Cells.Item($newRow, 11).Style = "Normal"

$ws.Cells.Item($newRow, 12).Value = "Phan Minh Khải"
$ws.Cells.Item($newRow, 13).Value = "5e552852d51d0258"
$ws.Cells.Item($newRow, 14).Value = "c20c40a4882271452afde009aad79f4bb64b82b8cb0a50a3552b1fdcc45af391"
